# Update row 20 (2025Q2) metrics on the active sheet:
#   total_customers (C20):     338 -> 344
#   returning_customers (D20): 265 -> 267
#   new_customers (E20):       73  -> 77
#   recurrence_rate (F20):     82.04334365325077 -> 82.6625386996904

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C20").Value = 344
$ws.Range("D20").Value = 267
$ws.Range("E20").Value = 77
$ws.Range("F20").Value = 82.6625386996904
